$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Tighten the trailing whitespace on the "vFPC (Lenny Colton)" line:
#    " (Lenny Colton)  " -> " (Lenny Colton) "
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "(Lenny Colton)  ", $true, $false, $false, $false, $false,
    $true, 1, $false, "(Lenny Colton) ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Locate the "vFPC (Lenny Colton)" list paragraph (by 1-based Item
#    index - the collection's .Index property is unreliable near the
#    end of the body, so always address paragraphs via Item(N)).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lennyNum = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Lenny Colton*") {
        $lennyNum = $i
    }
}

$lennyPara = $d.Paragraphs.Item($lennyNum)

# Insert a new paragraph right after it - Word clones the paragraph
# formatting (ListParagraph style, numPr, pBdr) of the paragraph it
# splits off of, which matches the new bullet's <w:pPr>.
$lennyPara.Range.InsertParagraphAfter()

$newParaNum = $lennyNum + 1
$newPara = $d.Paragraphs.Item($newParaNum)
$insAt = $newPara.Range.Start

# ---------------------------------------------------------------------
# 3. "VCH" as a hyperlink.
# ---------------------------------------------------------------------
$vchRange = $d.Range($insAt, $insAt)
$vchRange.InsertAfter("VCH")
$vchRange = $d.Range($insAt, $insAt + 3)
$d.Hyperlinks.Add($vchRange, "https://github.com/JanFries/VCH", [Type]::Missing, [Type]::Missing, "VCH") | Out-Null

# ---------------------------------------------------------------------
# 4. "(", "Jan Fries" and ") " as three separate plain-text runs,
#    appended right before the paragraph mark. The paragraph's End
#    keeps moving (the hyperlink field occupies hidden character slots)
#    so we re-query it before each insert.
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item($newParaNum)
$pos = $p.Range.End - 1
$d.Range($pos, $pos).InsertAfter("(") | Out-Null

$p = $d.Paragraphs.Item($newParaNum)
$pos = $p.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Jan Fries") | Out-Null

$p = $d.Paragraphs.Item($newParaNum)
$pos = $p.Range.End - 1
$d.Range($pos, $pos).InsertAfter(") ") | Out-Null
